$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "2,16 → 1,89 ↓↓↓", "↓↓↓ (over)")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "2,16 → 1,89 ↓↓↓", "↓↓↓ (over)")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "2,16 → 1,89 ↓↓↓", "↓↓↓ (over)")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "2,16 → 1,89 ↓↓↓", "↓↓↓ (over)")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "2,16 → 1,89 ↓↓↓", "↓↓↓ (over)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "1,63 → 1,90", "")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "1,63 → 1,90", "")
    ,@("Soccer, Club Friendlies", "Sparta Prague - FC Vysocina Jihlava", "Full time, (total) Under 4.0 / Over 4.0", "1,63 → 1,90", "")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "")
    ,@("Basketball, Australia - NBL1 Women", "Perth Redbacks - Warwick Senators", "Full time, Asian Handicap 34.5", "2,06 → 1,67 ↓↓↓", "")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Basketball, Australia - NBL1 Women", "Perth Redbacks - Warwick Senators", "Full time, Asian Handicap 34.5", "2,06 → 1,67 ↓↓↓", "↓↓↓ (home)")
    ,@("Basketball, Australia - NBL1 Women", "Perth Redbacks - Warwick Senators", "Full time, Asian Handicap 34.5", "2,06 → 1,67 ↓↓↓", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "1,49 → 1,97", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL New South Wales", "Sydney FC - Sydney United 58", "Full time, 3-way", "1,83 → 2,19", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL New South Wales", "Wollongong Wolves FC - St George City FA", "Full time, Asian Handicap 0.5", "1,62 → 1,77", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL New South Wales", "Wollongong Wolves FC - St George City FA", "Full time, Asian Handicap 0.5", "1,62 → 1,77", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL New South Wales", "Wollongong Wolves FC - St George City FA", "Full time, Asian Handicap 0.5", "2,30 → 1,91 ↓↓↓", "↓↓↓ (away)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "↓↓↓ (away)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "↓↓↓ (home)")
    ,@("Soccer, Australia - NPL South Australia Women", "Adelaide City - Fulham United", "Full time, Asian Handicap -2.5", "2,29 → 1,76 ↓↓↓", "↓↓↓ (home)")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = ""
}
